$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 18.00179481506348
$ws.Range("C3").Value = 17.69018173217773
$ws.Range("C4").Value = 17.26508140563965
$ws.Range("C5").Value = 19.06514167785645
$ws.Range("C6").Value = 18.15986633300781
